# Generate Report for Handback
# Replaces the two tracked source-file UUIDs (and their derived xliff file
# names / timestamps) across the "Overview", "zh-cn" and "de-de" sheets,
# mirroring a fresh handback-status report run.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "5861ba4d-4307-45f7-9a33-2f0b282d459c"
$oldUuid2 = "aa1a431c-2e36-47f3-b5b8-2b4d86ec199b"
$newUuid1 = "2340aa62-d89d-4552-a280-c1c763fff4c8"
$newUuid2 = "ffff6adbd3a3-6d85-4246-9e9b-10ec802cd590"

$newHash = "9e4496bbd74c383ea9b036cd6cb4ef053563b633"

$zhcnXlf = "$newUuid1.$newHash.zh-cn.xlf"
$dedeXlf = "$newUuid1.$newHash.de-de.xlf"

$latestHoDate = "2016-08-24 13:07:56"
$zhHandoffDate = "2016-08-24 13:07:51"
$zhHandbackDate = "2016-08-24 13:08:25"
$deHandbackDate = "2016-08-24 13:08:33"

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newUuid1.md"
$ws.Range("B2").Value = "e2e\$newUuid1.md"
$ws.Range("G2").Value = $latestHoDate

$ws.Range("A3").Value = "$newUuid2.md"
$ws.Range("B3").Value = "e2e\$newUuid2.md"
$ws.Range("G3").Value = $latestHoDate

$ws.Hyperlinks.Item(1).TextToDisplay = "e2e\$newUuid1.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "e2e\$newUuid2.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newUuid1.md"
$ws.Range("I2").Value = "$newUuid1.md"
$ws.Range("G2").Value = $zhcnXlf
$ws.Range("H2").Value = $zhHandoffDate
$ws.Range("J2").Value = $zhcnXlf
$ws.Range("K2").Value = $zhHandbackDate

$ws.Range("A3").Value = "$newUuid2.md"
$ws.Range("I3").Value = "$newUuid2.md"
$ws.Range("G3").Value = $zhcnXlf
$ws.Range("H3").Value = $zhHandoffDate
$ws.Range("J3").Value = $zhcnXlf
$ws.Range("K3").Value = $zhHandbackDate

$ws.Hyperlinks.Item(1).TextToDisplay = "$newUuid1.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "$newUuid1.md"
$ws.Hyperlinks.Item(3).TextToDisplay = "$newUuid2.md"
$ws.Hyperlinks.Item(4).TextToDisplay = "$newUuid2.md"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newUuid1.md"
$ws.Range("I2").Value = "$newUuid1.md"
$ws.Range("G2").Value = $dedeXlf
$ws.Range("H2").Value = $latestHoDate
$ws.Range("J2").Value = $dedeXlf
$ws.Range("K2").Value = $deHandbackDate

$ws.Range("A3").Value = "$newUuid2.md"
$ws.Range("I3").Value = "$newUuid2.md"
$ws.Range("G3").Value = $dedeXlf
$ws.Range("H3").Value = $latestHoDate
$ws.Range("J3").Value = $dedeXlf
$ws.Range("K3").Value = $deHandbackDate

$ws.Hyperlinks.Item(1).TextToDisplay = "$newUuid1.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "$newUuid1.md"
$ws.Hyperlinks.Item(3).TextToDisplay = "$newUuid2.md"
$ws.Hyperlinks.Item(4).TextToDisplay = "$newUuid2.md"
